$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35-58 down to 36-59.
$ws.Rows("35:35").Insert()

# Fill in the new row 35 with the data from the commit.
$ws.Range("A35").Value = 1
$ws.Range("B35").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C35").Value = "Arica y Parinacota"
$ws.Range("D35").Value = 44574
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = 100114001
$ws.Range("G35").Value = "Papa"
$ws.Range("H35").Value = "Asterix"
$ws.Range("I35").Value = "1a (cosecha)"
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 13000
$ws.Range("L35").Value = 14000
$ws.Range("M35").Value = 13500
$ws.Range("N35").Value = "$/malla 25 kilos"
$ws.Range("O35").Value = "Región del Maule"
$ws.Range("P35").Value = 540
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
